$d = $word.ActiveDocument

# 1. Text change: add "del gioco" before "cioè:"
$d.Content.Find.Execute(
    "eseguire tutte le mosse, cioè:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "eseguire tutte le mosse del gioco, cioè:",
    2
)

# 2. Text change: "verranno rimesse nel mazzo." -> "possono venir rimesse nel mazzo dall'utente."
$d.Content.Find.Execute(
    "non spostate in basi o sequenze verranno rimesse nel mazzo.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "non spostate in basi o sequenze possono venir rimesse nel mazzo dall’utente.",
    2
)

# 3. Text change: replace the "In caso di vittoria o resa..." paragraph text
$d.Content.Find.Execute(
    "In caso di vittoria o resa l’utente potrà scegliere se tornare al menù o fare un’altra partita.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "In caso di resa verrà chiesto all’utente se vuole riprovare o tornare al menù, mentre in caso di vincita potrà tornare menù.",
    2
)
